# Update CDA Logical model for ST.r2b
$wb = $excel.ActiveWorkbook

# 1. Rename the "Include from ParticipationTyp" sheet to "Include #0"
$includeWs = $wb.Worksheets.Item("Include from ParticipationTyp")
$includeWs.Name = "Include #0"

# 2. Update the Metadata sheet
$ws = $wb.Worksheets.Item("Metadata")

# Bump Version and Date values
$ws.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" property row right before "Description" (old row 11)
$ws.Rows.Item(11).Insert()

# Copy the formatting from the row above (Contact) so the new row matches
# the rest of the table's style, then overwrite with the new content.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

Write-Output "done"
